$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the title cell (A1): v36 -> v37 ---
$ws.Cells.Item(1, 1).Value = "sp_Blitz® Check ID List - v37 Nov 19, 2014"

# --- Add the two new "Wait Stats" check rows (217, 218) ---
$ws.Cells.Item(217, 1).Value = 152
$ws.Cells.Item(217, 2).Value = 240
$ws.Cells.Item(217, 3).Value = "Wait Stats"
$ws.Cells.Item(217, 4).Value = "Top Wait Stats"
$ws.Cells.Item(217, 5).Value = "http://BrentOzar.com/go/waits"
$ws.Hyperlinks.Add($ws.Cells.Item(217, 5), "http://BrentOzar.com/go/waits")
$ws.Cells.Item(217, 5).Style = "Hyperlink"

$ws.Cells.Item(218, 1).Value = 153
$ws.Cells.Item(218, 2).Value = 240
$ws.Cells.Item(218, 3).Value = "Wait Stats"
$ws.Cells.Item(218, 4).Value = "No Significant Waits Detected"
$ws.Cells.Item(218, 5).Value = "http://BrentOzar.com/go/waits"
$ws.Hyperlinks.Add($ws.Cells.Item(218, 5), "http://BrentOzar.com/go/waits")
$ws.Cells.Item(218, 5).Style = "Hyperlink"

# --- Move the selection down to the newly added row, like the author did ---
$ws.Range("A219").Select()
